$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "Metadata" sheet - update the StructureDefinition metadata values
# ---------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B3").Value  = "0.4.0-snapshot-1"                      # Version
$meta.Range("B6").Value  = "draft"                                  # Status
$meta.Range("B8").Value  = "2024-05-23T12:16:26+00:00"              # Date
$meta.Range("B10").Value = "ANS (https://esante.gouv.fr)"           # Contact

# ---------------------------------------------------------------------
# 2. "Elements" sheet - the two "Mapping" columns (AK = 37, AL = 38)
#    swapped places: "Mapping: RIM Mapping" moved from AK to AL, and
#    "Mapping: Spécification métier vers l'extension ROR DropZone"
#    moved from AL to AK - for every row (header + data).
# ---------------------------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

$lastRow = $elements.Cells.Item(1, 1).End(4).Row   # xlDown sentinel not needed; used below instead
$lastRow = 6

for ($r = 1; $r -le $lastRow; $r++) {
    $akCell = $elements.Cells.Item($r, 37)
    $alCell = $elements.Cells.Item($r, 38)

    $akVal = $akCell.Value2
    $alVal = $alCell.Value2

    $akCell.Value = $alVal
    $alCell.Value = $akVal
}

# Column widths follow the content: AK (was the narrow "RIM Mapping"
# column) becomes the wide column, AL becomes the narrow one.
$elements.Columns.Item(37).ColumnWidth = 64.35
$elements.Columns.Item(38).ColumnWidth = 24.18
